# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price/profit updates to the Zodiark_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 3383.611
$ws.Range("I32").Value = 2658.6
$ws.Range("J32").Value = 3662.4614
$ws.Range("K32").Value = 2658.6
$ws.Range("L32").Value = 3662.4614
$ws.Range("M32").Value = -2332.6
$ws.Range("N32").Value = -4314.4614

# Row 74
$ws.Range("H74").Value = 7650.3887
$ws.Range("I74").Value = 5813.5713
$ws.Range("K74").Value = 5813.5713
$ws.Range("M74").Value = -4877.5713

# Row 77
$ws.Range("H77").Value = 7650.3887
$ws.Range("I77").Value = 5813.5713
$ws.Range("K77").Value = 29067.8565
$ws.Range("M77").Value = -24387.8565

# Row 100
$ws.Range("H100").Value = 4169961
$ws.Range("I100").Value = 2889.4614
$ws.Range("J100").Value = 7356545
$ws.Range("K100").Value = 2889.4614
$ws.Range("L100").Value = 7356545
$ws.Range("M100").Value = -2348.4614
$ws.Range("N100").Value = -7357627

# Row 101
$ws.Range("H101").Value = 964.5454999999999
$ws.Range("I101").Value = 562
$ws.Range("J101").Value = 1300
$ws.Range("K101").Value = 1686
$ws.Range("L101").Value = 3900
$ws.Range("M101").Value = -64
$ws.Range("N101").Value = -7144

# Row 135
$ws.Range("H135").Value = 1705.5358
$ws.Range("I135").Value = 1721.3462
$ws.Range("K135").Value = 15492.1158
$ws.Range("M135").Value = -12957.1158

# Row 138
$ws.Range("H138").Value = 2548.1316
$ws.Range("I138").Value = 1186.4166
$ws.Range("J138").Value = 3176.6155
$ws.Range("K138").Value = 3559.2498
$ws.Range("L138").Value = 9529.8465
$ws.Range("M138").Value = 1580.7502
$ws.Range("N138").Value = -19809.8465

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2859.45
$ws.Range("I32").Value = 1446.2245
$ws.Range("J32").Value = 9154.727999999999
$ws.Range("K32").Value = 1446.2245
$ws.Range("L32").Value = 9154.727999999999
$ws.Range("M32").Value = -1159.2245
$ws.Range("N32").Value = -9728.727999999999

# Row 61
$ws.Range("H61").Value = 1879.9333
$ws.Range("I61").Value = 1624
$ws.Range("K61").Value = 1624
$ws.Range("M61").Value = -1412

# Row 132
$ws.Range("H132").Value = 7969.426
$ws.Range("I132").Value = 4238.82
$ws.Range("K132").Value = 12716.46
$ws.Range("M132").Value = -10186.46

# Row 136
$ws.Range("H136").Value = 1879.9333
$ws.Range("I136").Value = 1624
$ws.Range("K136").Value = 4872
$ws.Range("M136").Value = -2322

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2640.5881
$ws.Range("I20").Value = 2010.3334
$ws.Range("J20").Value = 2984.3635
$ws.Range("K20").Value = 2010.3334
$ws.Range("L20").Value = 2984.3635
$ws.Range("M20").Value = -1763.3334
$ws.Range("N20").Value = -3478.3635

# Row 134
$ws.Range("H134").Value = 3938.6775
$ws.Range("I134").Value = 4076.8845
$ws.Range("J134").Value = 3220
$ws.Range("K134").Value = 12230.6535
$ws.Range("L134").Value = 9660
$ws.Range("M134").Value = -9695.6535
$ws.Range("N134").Value = -14730

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1580.2727
$ws.Range("J22").Value = 2271.1428
$ws.Range("L22").Value = 2271.1428
$ws.Range("N22").Value = -2971.1428

# Row 31
$ws.Range("H31").Value = 4600
$ws.Range("I31").Value = 5900
$ws.Range("K31").Value = 5900
$ws.Range("M31").Value = -5605

# Row 34
$ws.Range("H34").Value = 4600
$ws.Range("I34").Value = 5900
$ws.Range("K34").Value = 5900
$ws.Range("M34").Value = -5698

# Row 42
$ws.Range("H42").Value = 6000
$ws.Range("J42").Value = 6000
$ws.Range("L42").Value = 6000
$ws.Range("N42").Value = -7186

# Row 99
$ws.Range("H99").Value = 3245.1482
$ws.Range("I99").Value = 3084.5
$ws.Range("J99").Value = 3418.1538
$ws.Range("K99").Value = 3084.5
$ws.Range("L99").Value = 3418.1538
$ws.Range("M99").Value = -1586.5
$ws.Range("N99").Value = -6414.1538

# Row 126
$ws.Range("H126").Value = 3245.1482
$ws.Range("I126").Value = 3084.5
$ws.Range("J126").Value = 3418.1538
$ws.Range("K126").Value = 9253.5
$ws.Range("L126").Value = 10254.4614
$ws.Range("M126").Value = -6783.5
$ws.Range("N126").Value = -15194.4614

# Row 134
$ws.Range("H134").Value = 818.85
$ws.Range("I134").Value = 809.8421
$ws.Range("K134").Value = 2429.5263
$ws.Range("M134").Value = 105.4737

$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 1300.52
$ws.Range("I132").Value = 737.6
$ws.Range("K132").Value = 6638.400000000001
$ws.Range("M132").Value = -4108.400000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 627
$ws.Range("I3").Value = 1347
$ws.Range("K3").Value = 1347
$ws.Range("M3").Value = -1231

# Row 29
$ws.Range("H29").Value = 29998.5
$ws.Range("I29").Value = 29997
$ws.Range("K29").Value = 29997
$ws.Range("M29").Value = -29707

# Row 123
$ws.Range("H123").Value = 75998
$ws.Range("J123").Value = 75998
$ws.Range("L123").Value = 75998
$ws.Range("N123").Value = -80898

$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 999.3333
$ws.Range("I4").Value = 999.3333
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 999.3333
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -886.3333
$ws.Range("N4").ClearContents()

# Row 28
$ws.Range("H28").Value = 999.3333
$ws.Range("I28").Value = 999.3333
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 999.3333
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -767.3333
$ws.Range("N28").ClearContents()

# Row 37
$ws.Range("H37").Value = 999.3333
$ws.Range("I37").Value = 999.3333
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 999.3333
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -892.3333
$ws.Range("N37").ClearContents()

# Row 46
$ws.Range("H46").Value = 1145.375
$ws.Range("I46").Value = 932.6
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 932.6
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -744.6
$ws.Range("N46").Value = -1876

# Row 100
$ws.Range("H100").Value = 3006.318
$ws.Range("J100").Value = 4892.75
$ws.Range("L100").Value = 4892.75
$ws.Range("N100").Value = -5974.75

# Row 136
$ws.Range("H136").Value = 3771.2188
$ws.Range("I136").Value = 3360.4443
$ws.Range("K136").Value = 10081.3329
$ws.Range("M136").Value = -7531.332900000001

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 710
$ws.Range("I29").Value = 710
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 710
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -420
$ws.Range("N29").ClearContents()
